$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "33.920.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.04%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.786.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.11%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "221.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.71%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.549"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  -4.02%  "

$ws.Range("E9").Value = "  +1.27%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0715"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.59%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0921"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.57%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.042.45"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.794.58"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.83%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.56%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.627"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.69%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "33.951.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.80%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.21"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.88%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.38%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0781"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.25%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.09%  "

$ws.Range("E22").Value = "  +2.67%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.47%  "

$ws.Range("E24").Value = "  -1.78%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.32%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.28%  "

$ws.Range("E28").Value = "  -2.11%  "

$ws.Range("E29").Value = "  +0.00%  "

$ws.Range("E30").Value = "  +0.75%  "

$ws.Range("E31").Value = "  +1.28%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.68"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.80%  "

$ws.Range("E33").Value = "  -2.25%  "

$ws.Range("E34").Value = "  -2.44%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.406.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.45%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.641"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.51%  "

$ws.Range("E37").Value = "  -0.16%  "

$ws.Range("E38").Value = "  -1.42%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.936"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.58%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "79.52"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.23%  "

$ws.Range("E41").Value = "  -3.25%  "

$ws.Range("E42").Value = "  -0.45%  "

$ws.Range("E43").Value = "  +1.84%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.94"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.44%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0494"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.30%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.941.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.32%  "

$ws.Range("E47").Value = "  -0.98%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "105.25"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.44%  "

$ws.Range("E49").Value = "  -0.20%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.81%  "

$ws.Range("E51").Value = "  -0.90%  "
